$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new column U with header "66accuracy"
$ws.Range("U1").Value = "66accuracy"

# Set U2:U11 to 1
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 21).Value = 1
}
